$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 112503698
$ws.Range("B6").Value = 96735
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."

# I6 holds a numeric-looking value but must stay text, like the source file.
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "100"
$ws.Range("I6").ClearFormats()

$ws.Range("P6").Value = "Gammal skog, Ög"
$ws.Range("Q6").Value = 564754
$ws.Range("R6").Value = 6511406
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Östergötland"
$ws.Range("U6").Value = "Norrköping"
$ws.Range("V6").Value = "Östergötland"
$ws.Range("W6").Value = "Simonstorp"

# Y6/AA6 hold date-looking text that must stay text (not become a real date serial).
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-10-01"
$ws.Range("Y6").ClearFormats()

$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-10-01"
$ws.Range("AA6").ClearFormats()

$ws.Range("AC6").Value = "Finns under en vält gran. Finns massvis i området runt omkring, hundratals."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Frida Blixt"
$ws.Range("AX6").Value = "Frida Blixt"

$wb.Save()
